# Update "opus_big Simple aWCE" sheet with re-run hyperparameter search results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("opus_big Simple aWCE")

# Validation BLEU (D), Runtime in Seconds (E), Compute in PFLOPs (G) for each
# hyperparameter combination row. Runtime (Hours) in column F is a formula
# (=E/3600) and recalculates automatically.
$ws.Range("D2").Value = 43.159799999999997
$ws.Range("E2").Value = 29016.542700000002
$ws.Range("G2").Value = 258.44920000000002

$ws.Range("D3").Value = 43.383899999999997
$ws.Range("E3").Value = 29011.771700000001
$ws.Range("G3").Value = 258.44920000000002

$ws.Range("D4").Value = 43.003300000000003
$ws.Range("E4").Value = 16505.728800000001
$ws.Range("G4").Value = 232.65119999999999

$ws.Range("D5").Value = 42.912300000000002
$ws.Range("E5").Value = 21470.682000000001
$ws.Range("G5").Value = 258.44920000000002

$ws.Range("D6").Value = 42.892600000000002
$ws.Range("E6").Value = 16748.4535
$ws.Range("G6").Value = 103.4242

$ws.Range("D7").Value = 42.933700000000002
$ws.Range("E7").Value = 17256.9139
$ws.Range("G7").Value = 232.65119999999999

$ws.Range("D8").Value = 43.305900000000001
$ws.Range("E8").Value = 27637.078300000001
$ws.Range("G8").Value = 258.44920000000002

$ws.Range("D9").Value = 42.904000000000003
$ws.Range("E9").Value = 21341.0075
$ws.Range("G9").Value = 258.44920000000002

$ws.Range("D10").Value = 42.998699999999999
$ws.Range("E10").Value = 16562.896799999999
$ws.Range("G10").Value = 103.4242

$ws.Range("D11").Value = 42.981699999999996
$ws.Range("E11").Value = 16649.787499999999
$ws.Range("G11").Value = 232.65119999999999

$ws.Range("D12").Value = 42.935099999999998
$ws.Range("E12").Value = 16950.433099999998
$ws.Range("G12").Value = 155.0676

$ws.Range("D13").Value = 42.993499999999997
$ws.Range("E13").Value = 25700.366300000002
$ws.Range("G13").Value = 232.65119999999999

$ws.Range("D14").Value = 43.179499999999997
$ws.Range("E14").Value = 27668.145400000001
$ws.Range("G14").Value = 258.44920000000002

$ws.Range("D15").Value = 43.644100000000002
$ws.Range("E15").Value = 27306.571
$ws.Range("G15").Value = 258.44920000000002

$ws.Range("D16").Value = 42.882100000000001
$ws.Range("E16").Value = 16452.034
$ws.Range("G16").Value = 232.65119999999999

$ws.Range("D17").Value = 43.441699999999997
$ws.Range("E17").Value = 19255.5524
$ws.Range("G17").Value = 129.2098

$ws.Range("D18").Value = 43.2774
$ws.Range("E18").Value = 31316.296999999999
$ws.Range("G18").Value = 310.13979999999998

$ws.Range("D19").Value = 43.0167
$ws.Range("E19").Value = 18583.752499999999
$ws.Range("G19").Value = 258.44920000000002

$ws.Range("D20").Value = 43.419800000000002
$ws.Range("E20").Value = 33099.654699999999
$ws.Range("G20").Value = 310.13979999999998

$ws.Range("D21").Value = 43.119500000000002
$ws.Range("E21").Value = 27529.222900000001
$ws.Range("G21").Value = 258.44920000000002

$ws.Range("D22").Value = 42.997799999999998
$ws.Range("E22").Value = 16602.3436
$ws.Range("G22").Value = 103.4242

$ws.Range("D23").Value = 42.849600000000002
$ws.Range("E23").Value = 17032.7451
$ws.Range("G23").Value = 103.4242

# Row 24 run failed / needs a rerun - flag it instead of filling in results.
$ws.Range("H24").Value = "RERUN"

$ws.Range("D25").Value = 43.133699999999997
$ws.Range("E25").Value = 29087.944100000001
$ws.Range("G25").Value = 258.44920000000002

$ws.Range("D26").Value = 42.942599999999999
$ws.Range("E26").Value = 16759.3397
$ws.Range("G26").Value = 103.4242

$ws.Range("D27").Value = 43.315899999999999
$ws.Range("E27").Value = 20323.4028
$ws.Range("G27").Value = 258.44920000000002

$ws.Range("D28").Value = 43.142600000000002
$ws.Range("E28").Value = 27060.105599999999
$ws.Range("G28").Value = 310.13979999999998

$ws.Range("D29").Value = 42.887999999999998
$ws.Range("E29").Value = 16448.672699999999
$ws.Range("G29").Value = 155.0676

$ws.Range("D30").Value = 43.148299999999999
$ws.Range("E30").Value = 20052.149799999999
$ws.Range("G30").Value = 258.44920000000002

$ws.Range("D31").Value = 42.9636
$ws.Range("E31").Value = 16479.568500000001
$ws.Range("G31").Value = 103.4242

$ws.Range("D32").Value = 42.554699999999997
$ws.Range("E32").Value = 12920.147499999999
$ws.Range("G32").Value = 103.4242

# Row 33 run also failed / needs a rerun - flag it instead of filling in results.
$ws.Range("H33").Value = "RERUN"

# This sheet is now the focus of attention - make it the active tab/selection.
$ws.Activate()
$ws.Range("I16").Select()
